$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.970.44'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '2.824.27'
$ws.Range('E3').Value = '  +3.33%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '354.11'
$ws.Range('E5').Value = '  +6.79%  '
$ws.Range('D6').Value = '113.82'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.600'
$ws.Range('E9').Value = '  +5.56%  '
$ws.Range('D10').Value = '42.01'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('D11').Value = '0.0851'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '20.08'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '7.70'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = '3.248.70'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').Value = '2.821.48'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = '0.897'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '51.871.06'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '7.34'
$ws.Range('E19').Value = '  +7.71%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '13.57'
$ws.Range('E21').Value = '  +1.70%  '
$ws.Range('D22').Value = '0.0₂01000'
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('D23').Value = '269.78'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('D24').Value = '69.65'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = '2.79'
$ws.Range('E25').Value = '  +5.96%  '
$ws.Range('D26').Value = '26.77'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '50.82'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('E32').Value = '  -3.16%  '
$ws.Range('D33').Value = '0.0451'
$ws.Range('E33').Value = '  +31.73%  '
$ws.Range('D34').Value = '5.85'
$ws.Range('E34').Value = '  +5.95%  '
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '2.09'
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '3.21'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '18.38'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').Value = '4.85'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('E41').Value = '  +7.33%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '23.67'
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '128.30'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '3.34'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').Value = '2.076.20'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  +4.34%  '
$ws.Range('D49').Value = '0.959'
$ws.Range('E49').Value = '  +10.61%  '
$ws.Range('E50').Value = '  +3.35%  '
$ws.Range('D51').Value = '60.43'
$ws.Range('E51').Value = '  +1.24%  '
